# zscore prop correct in classic ps
# Corrects the "Classical" column (column 2) values in the results table
# to reflect the fixed z-score proportion calculation.

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

$updates = @{
    2  = "-0.039    "
    3  = "(0.033)   "
    4  = "-0.046    "
    5  = "(0.039)   "
    6  = "0.041    "
    7  = "(0.043)   "
    8  = "-0.025    "
    9  = "(0.035)   "
    10 = "0.068    "
    11 = "(0.059)   "
    12 = "0.277 ***"
    13 = "(0.053)   "
    14 = "0.056    "
    15 = "(0.057)   "
    16 = "-0.056    "
    17 = "(0.035)   "
    18 = "0.040    "
    19 = "(0.039)   "
    20 = "0.016    "
    21 = "(0.038)   "
    22 = "-0.007    "
    23 = "(0.033)   "
    24 = "-0.002    "
    25 = "(0.035)   "
    26 = "0.138 ***"
    27 = "(0.038)   "
    28 = "-0.064    "
    29 = "(0.042)   "
    30 = "0.047    "
    31 = "(0.042)   "
    32 = "0.091    "
    33 = "(0.048)   "
    35 = "1.000    "
    36 = "0.625    "
    37 = "0.446    "
    38 = "0.839    "
    39 = "0.704    "
}

foreach ($rowIndex in $updates.Keys) {
    $cell = $table.Cell($rowIndex, 2)
    $cellRange = $cell.Range
    [void]$cellRange.MoveEnd(1, -1)
    $cellRange.Text = $updates[$rowIndex]
}
